$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 400  # H29: 220 -> 400
$ws.Cells.Item(29, 9).Value = 400  # I29: 220 -> 400
$ws.Cells.Item(29, 11).Value = 1200  # K29: 660 -> 1200
$ws.Cells.Item(29, 13).Value = -919  # M29: -379 -> -919

$ws.Cells.Item(39, 8).Value = 144.83333  # H39: 172.2 -> 144.83333
$ws.Cells.Item(39, 9).Value = 144.83333  # I39: 172.2 -> 144.83333
$ws.Cells.Item(39, 11).Value = 434.49999  # K39: 516.5999999999999 -> 434.49999
$ws.Cells.Item(39, 13).Value = -138.49999  # M39: -220.5999999999999 -> -138.49999

$ws.Cells.Item(80, 8).Value = 1927.75  # H80: 2180.5 -> 1927.75
$ws.Cells.Item(80, 9).Value = 129  # I80: 2000 -> 129
$ws.Cells.Item(80, 10).Value = 2527.3333  # J80: 2206.2856 -> 2527.3333
$ws.Cells.Item(80, 11).Value = 387  # K80: 6000 -> 387
$ws.Cells.Item(80, 12).Value = 7581.999899999999  # L80: 6618.8568 -> 7581.999899999999
$ws.Cells.Item(80, 13).Value = 611  # M80: -5002 -> 611
$ws.Cells.Item(80, 14).Value = -9577.999899999999  # N80: -8614.856800000001 -> -9577.999899999999

$ws.Cells.Item(83, 8).Value = 1927.75  # H83: 2180.5 -> 1927.75
$ws.Cells.Item(83, 9).Value = 129  # I83: 2000 -> 129
$ws.Cells.Item(83, 10).Value = 2527.3333  # J83: 2206.2856 -> 2527.3333
$ws.Cells.Item(83, 11).Value = 1161  # K83: 18000 -> 1161
$ws.Cells.Item(83, 12).Value = 22745.9997  # L83: 19856.5704 -> 22745.9997
$ws.Cells.Item(83, 13).Value = 3831  # M83: -13008 -> 3831
$ws.Cells.Item(83, 14).Value = -32729.9997  # N83: -29840.5704 -> -32729.9997

$ws.Cells.Item(106, 8).Value = 34000  # H106: 9500 -> 34000
$ws.Cells.Item(106, 9).Value = 34000  # I106: 0 -> 34000
$ws.Cells.Item(106, 10).Value = 0  # J106: 9500 -> 0
$ws.Cells.Item(106, 11).Value = 34000  # K106: 0 -> 34000
$ws.Cells.Item(106, 13).Value = -33369  # M106 added
$ws.Cells.Item(106, 14).Value = $null  # N106 removed (was -10762)

$ws.Cells.Item(112, 8).Value = 2596.1177  # H112: 2552.4443 -> 2596.1177
$ws.Cells.Item(112, 10).Value = 2596.1177  # J112: 2552.4443 -> 2596.1177
$ws.Cells.Item(112, 12).Value = 7788.353099999999  # L112: 7657.3329 -> 7788.353099999999
$ws.Cells.Item(112, 14).Value = -10004.3531  # N112: -9873.332900000001 -> -10004.3531

$ws.Cells.Item(127, 8).Value = 2917.25  # H127: 2433.1428 -> 2917.25
$ws.Cells.Item(127, 10).Value = 9000  # J127: 3590.75 -> 9000
$ws.Cells.Item(127, 12).Value = 27000  # L127: 10772.25 -> 27000
$ws.Cells.Item(127, 14).Value = -36920  # N127: -20692.25 -> -36920

$ws.Cells.Item(129, 8).Value = 2972.625  # H129: 3223.625 -> 2972.625
$ws.Cells.Item(129, 9).Value = 263  # I129: 266.33334 -> 263
$ws.Cells.Item(129, 10).Value = 4598.4  # J129: 4998 -> 4598.4
$ws.Cells.Item(129, 11).Value = 789  # K129: 799.0000200000001 -> 789
$ws.Cells.Item(129, 12).Value = 13795.2  # L129: 14994 -> 13795.2
$ws.Cells.Item(129, 13).Value = 4211  # M129: 4200.99998 -> 4211
$ws.Cells.Item(129, 14).Value = -23795.2  # N129: -24994 -> -23795.2

$ws.Cells.Item(132, 8).Value = 34486940  # H132: 34486964 -> 34486940
$ws.Cells.Item(132, 9).Value = 35718560  # I132: 35718584 -> 35718560
$ws.Cells.Item(132, 11).Value = 107155680  # K132: 107155752 -> 107155680
$ws.Cells.Item(132, 13).Value = -107153150  # M132: -107153222 -> -107153150

$ws.Cells.Item(138, 8).Value = 5343.7905  # H138: 5268.091 -> 5343.7905
$ws.Cells.Item(138, 9).Value = 5264.3076  # I138: 5031.0713 -> 5264.3076
$ws.Cells.Item(138, 10).Value = 5378.2334  # J138: 5378.7 -> 5378.2334
$ws.Cells.Item(138, 11).Value = 15792.9228  # K138: 15093.2139 -> 15792.9228
$ws.Cells.Item(138, 12).Value = 16134.7002  # L138: 16136.1 -> 16134.7002
$ws.Cells.Item(138, 13).Value = -10652.9228  # M138: -9953.213899999999 -> -10652.9228
$ws.Cells.Item(138, 14).Value = -26414.7002  # N138: -26416.1 -> -26414.7002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2798.25  # H2: 2766.5 -> 2798.25
$ws.Cells.Item(2, 9).Value = 2731  # I2: 2749.5 -> 2731
$ws.Cells.Item(2, 10).Value = 3000  # J2: 2775 -> 3000
$ws.Cells.Item(2, 11).Value = 2731  # K2: 2749.5 -> 2731
$ws.Cells.Item(2, 12).Value = 3000  # L2: 2775 -> 3000
$ws.Cells.Item(2, 13).Value = -2618  # M2: -2636.5 -> -2618
$ws.Cells.Item(2, 14).Value = -3226  # N2: -3001 -> -3226

$ws.Cells.Item(32, 8).Value = 18852.54  # H32: 17358.867 -> 18852.54
$ws.Cells.Item(32, 9).Value = 11518.5  # I32: 10873.75 -> 11518.5
$ws.Cells.Item(32, 11).Value = 11518.5  # K32: 10873.75 -> 11518.5
$ws.Cells.Item(32, 13).Value = -11231.5  # M32: -10586.75 -> -11231.5

$ws.Cells.Item(74, 8).Value = 5310490  # H74: 5752864.5 -> 5310490
$ws.Cells.Item(74, 9).Value = 6275126.5  # I74: 6902439.5 -> 6275126.5
$ws.Cells.Item(74, 11).Value = 6275126.5  # K74: 6902439.5 -> 6275126.5
$ws.Cells.Item(74, 13).Value = -6274252.5  # M74: -6901565.5 -> -6274252.5

$ws.Cells.Item(77, 8).Value = 5310490  # H77: 5752864.5 -> 5310490
$ws.Cells.Item(77, 9).Value = 6275126.5  # I77: 6902439.5 -> 6275126.5
$ws.Cells.Item(77, 11).Value = 31375632.5  # K77: 34512197.5 -> 31375632.5
$ws.Cells.Item(77, 13).Value = -31371264.5  # M77: -34507829.5 -> -31371264.5

$ws.Cells.Item(116, 8).Value = 2798.25  # H116: 2766.5 -> 2798.25
$ws.Cells.Item(116, 9).Value = 2731  # I116: 2749.5 -> 2731
$ws.Cells.Item(116, 10).Value = 3000  # J116: 2775 -> 3000
$ws.Cells.Item(116, 11).Value = 2731  # K116: 2749.5 -> 2731
$ws.Cells.Item(116, 12).Value = 3000  # L116: 2775 -> 3000
$ws.Cells.Item(116, 13).Value = -437  # M116: -455.5 -> -437
$ws.Cells.Item(116, 14).Value = -7588  # N116: -7363 -> -7588

$ws.Cells.Item(122, 8).Value = 4212.8237  # H122: 4674.1577 -> 4212.8237
$ws.Cells.Item(122, 9).Value = 3878.2  # I122: 4433.1763 -> 3878.2
$ws.Cells.Item(122, 11).Value = 11634.6  # K122: 13299.5289 -> 11634.6
$ws.Cells.Item(122, 13).Value = -9184.599999999999  # M122: -10849.5289 -> -9184.599999999999

$ws.Cells.Item(130, 8).Value = 47417.668  # H130: 48250.168 -> 47417.668
$ws.Cells.Item(130, 10).Value = 47417.668  # J130: 48250.168 -> 47417.668
$ws.Cells.Item(130, 12).Value = 47417.668  # L130: 48250.168 -> 47417.668
$ws.Cells.Item(130, 14).Value = -57457.668  # N130: -58290.168 -> -57457.668

$ws.Cells.Item(131, 8).Value = 79625  # H131: 79399 -> 79625
$ws.Cells.Item(131, 10).Value = 79625  # J131: 79399 -> 79625
$ws.Cells.Item(131, 12).Value = 79625  # L131: 79399 -> 79625
$ws.Cells.Item(131, 14).Value = -89705  # N131: -89479 -> -89705

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2798.25  # H3: 2766.5 -> 2798.25
$ws.Cells.Item(3, 9).Value = 2731  # I3: 2749.5 -> 2731
$ws.Cells.Item(3, 10).Value = 3000  # J3: 2775 -> 3000
$ws.Cells.Item(3, 11).Value = 2731  # K3: 2749.5 -> 2731
$ws.Cells.Item(3, 12).Value = 3000  # L3: 2775 -> 3000
$ws.Cells.Item(3, 13).Value = -2617  # M3: -2635.5 -> -2617
$ws.Cells.Item(3, 14).Value = -3228  # N3: -3003 -> -3228

$ws.Cells.Item(5, 8).Value = 247.625  # H5: 253.91667 -> 247.625
$ws.Cells.Item(5, 9).Value = 192.66667  # I5: 77.5 -> 192.66667
$ws.Cells.Item(5, 10).Value = 412.5  # J5: 606.75 -> 412.5
$ws.Cells.Item(5, 11).Value = 192.66667  # K5: 77.5 -> 192.66667
$ws.Cells.Item(5, 12).Value = 412.5  # L5: 606.75 -> 412.5
$ws.Cells.Item(5, 13).Value = -79.66667000000001  # M5: 35.5 -> -79.66667000000001
$ws.Cells.Item(5, 14).Value = -638.5  # N5: -832.75 -> -638.5

$ws.Cells.Item(97, 8).Value = 25350  # H97: 28000 -> 25350
$ws.Cells.Item(97, 9).Value = 25350  # I97: 28000 -> 25350
$ws.Cells.Item(97, 11).Value = 25350  # K97: 28000 -> 25350
$ws.Cells.Item(97, 13).Value = -24359  # M97: -27009 -> -24359

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3200  # H31: 2888.6667 -> 3200
$ws.Cells.Item(31, 9).Value = 0  # I31: 2332.6667 -> 0
$ws.Cells.Item(31, 10).Value = 3200  # J31: 3166.6667 -> 3200
$ws.Cells.Item(31, 11).Value = 0  # K31: 2332.6667 -> 0
$ws.Cells.Item(31, 13).Value = $null  # M31 removed (was -2037.6667)
$ws.Cells.Item(31, 14).Value = -3790  # N31: -3756.6667 -> -3790

$ws.Cells.Item(33, 8).Value = 1248  # H33: 1981.5 -> 1248
$ws.Cells.Item(33, 9).Value = 1248  # I33: 1515.5 -> 1248
$ws.Cells.Item(33, 10).Value = 0  # J33: 2447.5 -> 0
$ws.Cells.Item(33, 11).Value = 1248  # K33: 1515.5 -> 1248
$ws.Cells.Item(33, 12).Value = 0  # L33: 2447.5 -> 0
$ws.Cells.Item(33, 14).Value = $null  # N33 removed (was -3205.5)

$ws.Cells.Item(34, 8).Value = 3200  # H34: 2888.6667 -> 3200
$ws.Cells.Item(34, 9).Value = 0  # I34: 2332.6667 -> 0
$ws.Cells.Item(34, 10).Value = 3200  # J34: 3166.6667 -> 3200
$ws.Cells.Item(34, 11).Value = 0  # K34: 2332.6667 -> 0
$ws.Cells.Item(34, 13).Value = $null  # M34 removed (was -2130.6667)
$ws.Cells.Item(34, 14).Value = -3604  # N34: -3570.6667 -> -3604

$ws.Cells.Item(105, 8).Value = 1698.3334  # H105: 1700 -> 1698.3334
$ws.Cells.Item(105, 9).Value = 1698.3334  # I105: 1700 -> 1698.3334
$ws.Cells.Item(105, 11).Value = 1698.3334  # K105: 1700 -> 1698.3334
$ws.Cells.Item(105, 13).Value = 48.66660000000002  # M105: 47 -> 48.66660000000002

$ws.Cells.Item(109, 8).Value = 35000  # H109: 0 -> 35000
$ws.Cells.Item(109, 10).Value = 35000  # J109: 0 -> 35000
$ws.Cells.Item(109, 12).Value = 35000  # L109: 0 -> 35000
$ws.Cells.Item(109, 14).Value = -37080  # N109 added

$ws.Cells.Item(129, 8).Value = 32122.5  # H129: 37396 -> 32122.5
$ws.Cells.Item(129, 9).Value = 26248.5  # I129: 34994 -> 26248.5
$ws.Cells.Item(129, 11).Value = 26248.5  # K129: 34994 -> 26248.5
$ws.Cells.Item(129, 13).Value = -21248.5  # M129: -29994 -> -21248.5

$ws.Cells.Item(132, 8).Value = 15395429  # H132: 14295886 -> 15395429
$ws.Cells.Item(132, 9).Value = 16678131  # I132: 15395339 -> 16678131
$ws.Cells.Item(132, 10).Value = 3000  # J132: 2999 -> 3000
$ws.Cells.Item(132, 11).Value = 50034393  # K132: 46186017 -> 50034393
$ws.Cells.Item(132, 12).Value = 9000  # L132: 8997 -> 9000
$ws.Cells.Item(132, 13).Value = -50031863  # M132: -46183487 -> -50031863
$ws.Cells.Item(132, 14).Value = -14060  # N132: -14057 -> -14060

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value = 515.7143  # H6: 446.92856 -> 515.7143
$ws.Cells.Item(6, 9).Value = 185  # I6: 173.61539 -> 185
$ws.Cells.Item(6, 10).Value = 2500  # J6: 4000 -> 2500
$ws.Cells.Item(6, 11).Value = 555  # K6: 520.84617 -> 555
$ws.Cells.Item(6, 12).Value = 7500  # L6: 12000 -> 7500
$ws.Cells.Item(6, 13).Value = -442  # M6: -407.84617 -> -442
$ws.Cells.Item(6, 14).Value = -7726  # N6: -12226 -> -7726

$ws.Cells.Item(44, 8).Value = 125243.375  # H44: 111331.22 -> 125243.375
$ws.Cells.Item(44, 9).Value = 200041.4  # I44: 166706.83 -> 200041.4
$ws.Cells.Item(44, 11).Value = 600124.2  # K44: 500120.49 -> 600124.2
$ws.Cells.Item(44, 13).Value = -599726.2  # M44: -499722.49 -> -599726.2

$ws.Cells.Item(69, 8).Value = 2665.9524  # H69: 2662.093 -> 2665.9524
$ws.Cells.Item(69, 10).Value = 2607.7568  # J69: 2604.9211 -> 2607.7568
$ws.Cells.Item(69, 12).Value = 7823.2704  # L69: 7814.763300000001 -> 7823.2704
$ws.Cells.Item(69, 14).Value = -9445.270400000001  # N69: -9436.763300000001 -> -9445.270400000001

$ws.Cells.Item(72, 8).Value = 2665.9524  # H72: 2662.093 -> 2665.9524
$ws.Cells.Item(72, 10).Value = 2607.7568  # J72: 2604.9211 -> 2607.7568
$ws.Cells.Item(72, 12).Value = 23469.8112  # L72: 23444.2899 -> 23469.8112
$ws.Cells.Item(72, 14).Value = -31581.8112  # N72: -31556.2899 -> -31581.8112

$ws.Cells.Item(129, 8).Value = 2639.9  # H129: 2721.4 -> 2639.9
$ws.Cells.Item(129, 9).Value = 1919.25  # I129: 2097.25 -> 1919.25
$ws.Cells.Item(129, 10).Value = 3120.3333  # J129: 3137.5 -> 3120.3333
$ws.Cells.Item(129, 11).Value = 5757.75  # K129: 6291.75 -> 5757.75
$ws.Cells.Item(129, 12).Value = 9360.999899999999  # L129: 9412.5 -> 9360.999899999999
$ws.Cells.Item(129, 13).Value = -757.75  # M129: -1291.75 -> -757.75
$ws.Cells.Item(129, 14).Value = -19360.9999  # N129: -19412.5 -> -19360.9999

$ws.Cells.Item(131, 8).Value = 1497  # H131: 1113.2222 -> 1497
$ws.Cells.Item(131, 9).Value = 1494  # I131: 941.3333 -> 1494
$ws.Cells.Item(131, 10).Value = 1500  # J131: 1199.1666 -> 1500
$ws.Cells.Item(131, 11).Value = 4482  # K131: 2823.9999 -> 4482
$ws.Cells.Item(131, 12).Value = 4500  # L131: 3597.4998 -> 4500
$ws.Cells.Item(131, 13).Value = 558  # M131: 2216.0001 -> 558
$ws.Cells.Item(131, 14).Value = -14580  # N131: -13677.4998 -> -14580

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 2552.0715  # H97: 2398.7778 -> 2552.0715
$ws.Cells.Item(97, 9).Value = 2543.125  # I97: 2231 -> 2543.125
$ws.Cells.Item(97, 10).Value = 2564  # J97: 2608.5 -> 2564
$ws.Cells.Item(97, 11).Value = 2543.125  # K97: 2231 -> 2543.125
$ws.Cells.Item(97, 12).Value = 2564  # L97: 2608.5 -> 2564
$ws.Cells.Item(97, 13).Value = -2047.125  # M97: -1735 -> -2047.125
$ws.Cells.Item(97, 14).Value = -3556  # N97: -3600.5 -> -3556

$ws.Cells.Item(102, 8).Value = 2605.1765  # H102: 2594.0527 -> 2605.1765
$ws.Cells.Item(102, 9).Value = 2654  # I102: 2631.9285 -> 2654
$ws.Cells.Item(102, 11).Value = 2654  # K102: 2631.9285 -> 2654
$ws.Cells.Item(102, 13).Value = -1032  # M102: -1009.9285 -> -1032

$ws.Cells.Item(122, 8).Value = 4720.2  # H122: 4826.5264 -> 4720.2
$ws.Cells.Item(122, 9).Value = 4650.923  # I122: 4813.5 -> 4650.923
$ws.Cells.Item(122, 11).Value = 13952.769  # K122: 14440.5 -> 13952.769
$ws.Cells.Item(122, 13).Value = -11502.769  # M122: -11990.5 -> -11502.769

$ws.Cells.Item(128, 8).Value = 97268  # H128: 93329.336 -> 97268
$ws.Cells.Item(128, 10).Value = 97268  # J128: 93329.336 -> 97268
$ws.Cells.Item(128, 12).Value = 97268  # L128: 93329.336 -> 97268
$ws.Cells.Item(128, 14).Value = -107228  # N128: -103289.336 -> -107228

$ws.Cells.Item(132, 8).Value = 23811888  # H132: 25643418 -> 23811888
$ws.Cells.Item(132, 9).Value = 2548.5833  # I132: 2598.9092 -> 2548.5833
$ws.Cells.Item(132, 11).Value = 7645.749899999999  # K132: 7796.7276 -> 7645.749899999999
$ws.Cells.Item(132, 13).Value = -5115.749899999999  # M132: -5266.7276 -> -5115.749899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4539.6  # H7: 5199.5 -> 4539.6
$ws.Cells.Item(7, 9).Value = 5174.5  # I7: 6266 -> 5174.5
$ws.Cells.Item(7, 11).Value = 5174.5  # K7: 6266 -> 5174.5
$ws.Cells.Item(7, 13).Value = -5062.5  # M7: -6154 -> -5062.5

$ws.Cells.Item(46, 8).Value = 500  # H46: 1244.5 -> 500
$ws.Cells.Item(46, 9).Value = 500  # I46: 1190 -> 500
$ws.Cells.Item(46, 10).Value = 0  # J46: 1299 -> 0
$ws.Cells.Item(46, 11).Value = 500  # K46: 1190 -> 500
$ws.Cells.Item(46, 12).Value = 0  # L46: 1299 -> 0
$ws.Cells.Item(46, 14).Value = $null  # N46 removed (was -1675)

$ws.Cells.Item(126, 8).Value = 4539.6  # H126: 5199.5 -> 4539.6
$ws.Cells.Item(126, 9).Value = 5174.5  # I126: 6266 -> 5174.5
$ws.Cells.Item(126, 11).Value = 15523.5  # K126: 18798 -> 15523.5
$ws.Cells.Item(126, 13).Value = -13053.5  # M126: -16328 -> -13053.5

$ws.Cells.Item(130, 8).Value = 76276.336  # H130: 87212.25 -> 76276.336
$ws.Cells.Item(130, 10).Value = 76276.336  # J130: 87212.25 -> 76276.336
$ws.Cells.Item(130, 12).Value = 76276.336  # L130: 87212.25 -> 76276.336
$ws.Cells.Item(130, 14).Value = -86316.336  # N130: -97252.25 -> -86316.336

$ws.Cells.Item(132, 8).Value = 8303.23  # H132: 8295.571 -> 8303.23
$ws.Cells.Item(132, 9).Value = 5168.875  # I132: 5505.222 -> 5168.875
$ws.Cells.Item(132, 11).Value = 15506.625  # K132: 16515.666 -> 15506.625
$ws.Cells.Item(132, 13).Value = -12976.625  # M132: -13985.666 -> -12976.625

$ws.Cells.Item(136, 8).Value = 2122.8096  # H136: 2154.2 -> 2122.8096
$ws.Cells.Item(136, 9).Value = 2088.3684  # I136: 2121.3333 -> 2088.3684
$ws.Cells.Item(136, 11).Value = 6265.1052  # K136: 6363.999899999999 -> 6265.1052
$ws.Cells.Item(136, 13).Value = -3715.1052  # M136: -3813.999899999999 -> -3715.1052

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 853.44446  # H113: 974 -> 853.44446
$ws.Cells.Item(113, 9).Value = 853.44446  # I113: 870.5714 -> 853.44446
$ws.Cells.Item(113, 10).Value = 0  # J113: 1698 -> 0
$ws.Cells.Item(113, 11).Value = 2560.33338  # K113: 2611.7142 -> 2560.33338
$ws.Cells.Item(113, 12).Value = 0  # L113: 5094 -> 0
$ws.Cells.Item(113, 14).Value = $null  # N113 removed (was -9434)

$ws.Cells.Item(122, 8).Value = 3264.8462  # H122: 3120.25 -> 3264.8462
$ws.Cells.Item(122, 10).Value = 5000  # J122: 0 -> 5000
$ws.Cells.Item(122, 12).Value = 15000  # L122: 0 -> 15000
$ws.Cells.Item(122, 14).Value = -19900  # N122 added

$ws.Cells.Item(126, 8).Value = 6945503.5  # H126: 8334351 -> 6945503.5
$ws.Cells.Item(126, 9).Value = 8334350.5  # I126: 10417622 -> 8334350.5
$ws.Cells.Item(126, 11).Value = 25003051.5  # K126: 31252866 -> 25003051.5
$ws.Cells.Item(126, 13).Value = -25000581.5  # M126: -31250396 -> -25000581.5

$ws.Cells.Item(132, 8).Value = 333334000  # H132: 1000000000 -> 333334000
$ws.Cells.Item(132, 10).Value = 333334000  # J132: 1000000000 -> 333334000
$ws.Cells.Item(132, 12).Value = 1000002000  # L132: 3000000000 -> 1000002000
$ws.Cells.Item(132, 14).Value = -1000007060  # N132: -3000005060 -> -1000007060

$ws.Cells.Item(136, 8).Value = 6658.0415  # H136: 6903.9565 -> 6658.0415
$ws.Cells.Item(136, 9).Value = 6861.591  # I136: 7140.619 -> 6861.591
$ws.Cells.Item(136, 11).Value = 20584.773  # K136: 21421.857 -> 20584.773
$ws.Cells.Item(136, 13).Value = -18034.773  # M136: -18871.857 -> -18034.773
